$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Day11 solution runtimes entered into row 15 (B = part1 time, C = part2 time)
$ws.Range("B15").Value = 0.36863070004619602
$ws.Range("C15").Value = 0.00064990000100806301

# Recalculate so the dependent SUM/shared formulas (E15, B31, C31, E31) update
$excel.Calculate()

# Update the selected/active cell to reflect where the user left off
$ws.Range("E15").Select()
